$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.646.68'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.64%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.549.48'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.94%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '504.02'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.60'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -8.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.549'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -6.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.556.19'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.17'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -7.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.100'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.328'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -5.61%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.002.66'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.625.66'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.47'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -5.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000133'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -5.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.555.37'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.48'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -5.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '329.65'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -7.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.99'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -5.56%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.90'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -5.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '59.15'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.404'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -5.04%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.157'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -6.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0769'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -9.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.82'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -7.78%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '149.21'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.45'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.77'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -7.87%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.91'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -6.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.877'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.97%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -8.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '35.65'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.820'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -9.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '284.40'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.37'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -8.64%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.46'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -8.96%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0975'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.73%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.602'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0527'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -5.87%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.45'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.99%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.49'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -8.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.901.73'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.02%  '
